$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = 48.091872
$ws.Range("H2").Value = 144.275616
$ws.Range("I2").Value = 0.421093842675958
$ws.Range("J2").Value = 0.423782205092405
$ws.Range("M2").Value = 121.928739
$ws.Range("N2").Value = 365.786217
$ws.Range("O2").Value = 0.2282232151508951
$ws.Range("P2").Value = 0.2419720431319445
$ws.Range("Q2").Value = 5863.781309109408
$ws.Range("R2").Value = 52774.03178198468
$ws.Range("S2").Value = 0.09610339065575232
$ws.Range("T2").Value = 0.10254344600917
$ws.Range("G3").Value = 48.091872
$ws.Range("H3").Value = 144.275616
$ws.Range("I3").Value = 0.421093842675958
$ws.Range("J3").Value = 0.423782205092405
$ws.Range("O3").Value = 0.2768624053389947
$ws.Range("P3").Value = 0.2935413991166814
$ws.Range("Q3").Value = 7113.477025325759
$ws.Range("R3").Value = 64021.29322793184
$ws.Range("S3").Value = 0.1165850541567059
$ws.Range("T3").Value = 0.124397621403577
$ws.Range("G4").Value = 48.091872
$ws.Range("H4").Value = 144.275616
$ws.Range("I4").Value = 0.421093842675958
$ws.Range("J4").Value = 0.423782205092405
$ws.Range("M4").Value = 83.50496933333334
$ws.Range("N4").Value = 250.514908
$ws.Range("O4").Value = 0.1563025480180701
$ws.Range("P4").Value = 0.1657186665504434
$ws.Range("Q4").Value = 4015.910296542592
$ws.Range("R4").Value = 36143.19266888333
$ws.Range("S4").Value = 0.06581804056497256
$ws.Range("T4").Value = 0.07022862193571987
$ws.Range("G5").Value = 48.091872
$ws.Range("H5").Value = 144.275616
$ws.Range("I5").Value = 0.421093842675958
$ws.Range("J5").Value = 0.423782205092405
$ws.Range("M5").Value = 91.06846250000001
$ws.Range("N5").Value = 182.136925
$ws.Range("O5").Value = 0.1704597085236707
$ws.Range("P5").Value = 0.1204857969594293
$ws.Range("Q5").Value = 4379.6528417868
$ws.Range("R5").Value = 26277.9170507208
$ws.Range("S5").Value = 0.07177953368365626
$ws.Range("T5").Value = 0.05105973671778275
$ws.Range("G6").Value = 48.091872
$ws.Range("H6").Value = 144.275616
$ws.Range("I6").Value = 0.421093842675958
$ws.Range("J6").Value = 0.423782205092405
$ws.Range("M6").Value = 89.83562999999999
$ws.Range("N6").Value = 269.50689
$ws.Range("O6").Value = 0.1681521229683693
$ws.Range("P6").Value = 0.1782820942415013
$ws.Range("Q6").Value = 4320.36361899936
$ws.Range("R6").Value = 38883.27257099425
$ws.Range("S6").Value = 0.07080782361487084
$ws.Range("T6").Value = 0.07555277902615537
$ws.Range("I7").Value = 0.1230362686979479
$ws.Range("J7").Value = 0.1238217612582891
$ws.Range("M7").Value = 121.928739
$ws.Range("N7").Value = 365.786217
$ws.Range("O7").Value = 0.2282232151508951
$ws.Range("P7").Value = 0.2419720431319445
$ws.Range("Q7").Value = 1713.294519219008
$ws.Range("R7").Value = 15419.65067297107
$ws.Range("S7").Value = 0.02807973282241511
$ws.Range("T7").Value = 0.02996140455586406
$ws.Range("I8").Value = 0.1230362686979479
$ws.Range("J8").Value = 0.1238217612582891
$ws.Range("O8").Value = 0.2768624053389947
$ws.Range("P8").Value = 0.2935413991166814
$ws.Range("Q8").Value = 2078.433788304426
$ws.Range("S8").Value = 0.03406411729564873
$ws.Range("T8").Value = 0.03634681304084987
$ws.Range("I9").Value = 0.1230362686979479
$ws.Range("J9").Value = 0.1238217612582891
$ws.Range("M9").Value = 83.50496933333334
$ws.Range("N9").Value = 250.514908
$ws.Range("O9").Value = 0.1563025480180701
$ws.Range("P9").Value = 0.1657186665504434
$ws.Range("Q9").Value = 1173.378872444103
$ws.Range("R9").Value = 10560.40985199693
$ws.Range("S9").Value = 0.01923088229612518
$ws.Range("T9").Value = 0.02051957716565101
$ws.Range("I10").Value = 0.1230362686979479
$ws.Range("J10").Value = 0.1238217612582891
$ws.Range("M10").Value = 91.06846250000001
$ws.Range("N10").Value = 182.136925
$ws.Range("O10").Value = 0.1704597085236707
$ws.Range("P10").Value = 0.1204857969594293
$ws.Range("Q10").Value = 1279.658093363467
$ws.Range("R10").Value = 7677.9485601808
$ws.Range("S10").Value = 0.02097272650009224
$ws.Range("T10").Value = 0.01491876358612515
$ws.Range("I11").Value = 0.1230362686979479
$ws.Range("J11").Value = 0.1238217612582891
$ws.Range("M11").Value = 89.83562999999999
$ws.Range("N11").Value = 269.50689
$ws.Range("O11").Value = 0.1681521229683693
$ws.Range("P11").Value = 0.1782820942415013
$ws.Range("Q11").Value = 1262.33481763136
$ws.Range("R11").Value = 11361.01335868224
$ws.Range("S11").Value = 0.02068880978366667
$ws.Range("T11").Value = 0.02207520290979897
$ws.Range("G12").Value = 21.412221
$ws.Range("H12").Value = 64.23666299999999
$ws.Range("I12").Value = 0.1874860355013181
$ws.Range("J12").Value = 0.1886829905749125
$ws.Range("M12").Value = 121.928739
$ws.Range("N12").Value = 365.786217
$ws.Range("O12").Value = 0.2282232151508951
$ws.Range("P12").Value = 0.2419720431319445
$ws.Range("Q12").Value = 2610.765105719319
$ws.Range("R12").Value = 23496.88595147387
$ws.Range("S12").Value = 0.04278866581800567
$ws.Range("T12").Value = 0.04565600873365702
$ws.Range("G13").Value = 21.412221
$ws.Range("H13").Value = 64.23666299999999
$ws.Range("I13").Value = 0.1874860355013181
$ws.Range("J13").Value = 0.1886829905749125
$ws.Range("O13").Value = 0.2768624053389947
$ws.Range("P13").Value = 0.2935413991166814
$ws.Range("Q13").Value = 3167.174323026929
$ws.Range("R13").Value = 28504.56890724236
$ws.Range("S13").Value = 0.05190783475636707
$ws.Range("T13").Value = 0.05538626904287943
$ws.Range("G14").Value = 21.412221
$ws.Range("H14").Value = 64.23666299999999
$ws.Range("I14").Value = 0.1874860355013181
$ws.Range("J14").Value = 0.1886829905749125
$ws.Range("M14").Value = 83.50496933333334
$ws.Range("N14").Value = 250.514908
$ws.Range("O14").Value = 0.1563025480180701
$ws.Range("P14").Value = 0.1657186665504434
$ws.Range("Q14").Value = 1788.026857963556
$ws.Range("R14").Value = 16092.241721672
$ws.Range("S14").Value = 0.02930454506666235
$ws.Range("T14").Value = 0.03126829359882437
$ws.Range("G15").Value = 21.412221
$ws.Range("H15").Value = 64.23666299999999
$ws.Range("I15").Value = 0.1874860355013181
$ws.Range("J15").Value = 0.1886829905749125
$ws.Range("M15").Value = 91.06846250000001
$ws.Range("N15").Value = 182.136925
$ws.Range("O15").Value = 0.1704597085236707
$ws.Range("P15").Value = 0.1204857969594293
$ws.Range("Q15").Value = 1949.978045180213
$ws.Range("R15").Value = 11699.86827108128
$ws.Range("S15").Value = 0.03195881496381326
$ws.Range("T15").Value = 0.02273362049210683
$ws.Range("G16").Value = 21.412221
$ws.Range("H16").Value = 64.23666299999999
$ws.Range("I16").Value = 0.1874860355013181
$ws.Range("J16").Value = 0.1886829905749125
$ws.Range("M16").Value = 89.83562999999999
$ws.Range("N16").Value = 269.50689
$ws.Range("O16").Value = 0.1681521229683693
$ws.Range("P16").Value = 0.1782820942415013
$ws.Range("Q16").Value = 1923.58036323423
$ws.Range("R16").Value = 17312.22326910807
$ws.Range("S16").Value = 0.03152617489646969
$ws.Range("T16").Value = 0.03363879870744485
$ws.Range("G17").Value = 2.1734975
$ws.Range("H17").Value = 4.346995
$ws.Range("I17").Value = 0.01903120789977957
$ws.Range("J17").Value = 0.012768471746644
$ws.Range("M17").Value = 121.928739
$ws.Range("N17").Value = 365.786217
$ws.Range("O17").Value = 0.2282232151508951
$ws.Range("P17").Value = 0.2419720431319445
$ws.Range("Q17").Value = 265.0118093946525
$ws.Range("R17").Value = 1590.070856367915
$ws.Range("S17").Value = 0.004343363455092807
$ws.Range("T17").Value = 0.003089613196207957
$ws.Range("G18").Value = 2.1734975
$ws.Range("H18").Value = 4.346995
$ws.Range("I18").Value = 0.01903120789977957
$ws.Range("J18").Value = 0.012768471746644
$ws.Range("O18").Value = 0.2768624053389947
$ws.Range("P18").Value = 0.2935413991166814
$ws.Range("Q18").Value = 321.4914264691749
$ws.Range("R18").Value = 1928.94855881505
$ws.Range("S18").Value = 0.005269025995639449
$ws.Range("T18").Value = 0.003748075061091696
$ws.Range("G19").Value = 2.1734975
$ws.Range("H19").Value = 4.346995
$ws.Range("I19").Value = 0.01903120789977957
$ws.Range("J19").Value = 0.012768471746644
$ws.Range("M19").Value = 83.50496933333334
$ws.Range("N19").Value = 250.514908
$ws.Range("O19").Value = 0.1563025480180701
$ws.Range("P19").Value = 0.1657186665504434
$ws.Range("Q19").Value = 181.4978420835767
$ws.Range("R19").Value = 1088.98705250146
$ws.Range("S19").Value = 0.002974626286597171
$ws.Range("T19").Value = 0.002115974111740854
$ws.Range("G20").Value = 2.1734975
$ws.Range("H20").Value = 4.346995
$ws.Range("I20").Value = 0.01903120789977957
$ws.Range("J20").Value = 0.012768471746644
$ws.Range("M20").Value = 91.06846250000001
$ws.Range("N20").Value = 182.136925
$ws.Range("O20").Value = 0.1704597085236707
$ws.Range("P20").Value = 0.1204857969594293
$ws.Range("Q20").Value = 197.9370755725938
$ws.Range("R20").Value = 791.748302290375
$ws.Range("S20").Value = 0.003244054151449806
$ws.Range("T20").Value = 0.001538419494348359
$ws.Range("G21").Value = 2.1734975
$ws.Range("H21").Value = 4.346995
$ws.Range("I21").Value = 0.01903120789977957
$ws.Range("J21").Value = 0.012768471746644
$ws.Range("M21").Value = 89.83562999999999
$ws.Range("N21").Value = 269.50689
$ws.Range("O21").Value = 0.1681521229683693
$ws.Range("P21").Value = 0.1782820942415013
$ws.Range("Q21").Value = 195.257517215925
$ws.Range("R21").Value = 1171.54510329555
$ws.Range("S21").Value = 0.003200138011000336
$ws.Range("T21").Value = 0.002276389883255132
$ws.Range("G22").Value = 28.477822
$ws.Range("H22").Value = 85.433466
$ws.Range("I22").Value = 0.2493526452249964
$ws.Range("J22").Value = 0.2509445713277496
$ws.Range("M22").Value = 121.928739
$ws.Range("N22").Value = 365.786217
$ws.Range("O22").Value = 0.2282232151508951
$ws.Range("P22").Value = 0.2419720431319445
$ws.Range("Q22").Value = 3472.264925926458
$ws.Range("R22").Value = 31250.38433333812
$ws.Range("S22").Value = 0.05690806239962916
$ws.Range("T22").Value = 0.06072157063704555
$ws.Range("G23").Value = 28.477822
$ws.Range("H23").Value = 85.433466
$ws.Range("I23").Value = 0.2493526452249964
$ws.Range("J23").Value = 0.2509445713277496
$ws.Range("O23").Value = 0.2768624053389947
$ws.Range("P23").Value = 0.2935413991166814
$ws.Range("Q23").Value = 4212.27796098926
$ws.Range("R23").Value = 37910.50164890334
$ws.Range("S23").Value = 0.06903637313463348
$ws.Range("T23").Value = 0.07366262056828346
$ws.Range("G24").Value = 28.477822
$ws.Range("H24").Value = 85.433466
$ws.Range("I24").Value = 0.2493526452249964
$ws.Range("J24").Value = 0.2509445713277496
$ws.Range("M24").Value = 83.50496933333334
$ws.Range("N24").Value = 250.514908
$ws.Range("O24").Value = 0.1563025480180701
$ws.Range("P24").Value = 0.1657186665504434
$ws.Range("Q24").Value = 2378.039652790125
$ws.Range("R24").Value = 21402.35687511113
$ws.Range("S24").Value = 0.03897445380371278
$ws.Range("T24").Value = 0.04158619973850728
$ws.Range("G25").Value = 28.477822
$ws.Range("H25").Value = 85.433466
$ws.Range("I25").Value = 0.2493526452249964
$ws.Range("J25").Value = 0.2509445713277496
$ws.Range("M25").Value = 91.06846250000001
$ws.Range("N25").Value = 182.136925
$ws.Range("O25").Value = 0.1704597085236707
$ws.Range("P25").Value = 0.1204857969594293
$ws.Range("Q25").Value = 2593.431464888675
$ws.Range("R25").Value = 15560.58878933205
$ws.Range("S25").Value = 0.04250457922465916
$ws.Range("T25").Value = 0.03023525666906627
$ws.Range("G26").Value = 28.477822
$ws.Range("H26").Value = 85.433466
$ws.Range("I26").Value = 0.2493526452249964
$ws.Range("J26").Value = 0.2509445713277496
$ws.Range("M26").Value = 89.83562999999999
$ws.Range("N26").Value = 269.50689
$ws.Range("O26").Value = 0.1681521229683693
$ws.Range("P26").Value = 0.1782820942415013
$ws.Range("Q26").Value = 2558.32308039786
$ws.Range("R26").Value = 23024.90772358074
$ws.Range("S26").Value = 0.04192917666236175
$ws.Range("T26").Value = 0.04473892371484699
